$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.2096317280453258
$ws.Cells.Item(2, 3).Value = 0.5184135977337111
$ws.Cells.Item(2, 10).Value = 0.028328611898017
$ws.Cells.Item(2, 16).Value = 0.1359773371104816
$ws.Cells.Item(2, 19).Value = 0.1076487252124646
$ws.Cells.Item(3, 2).Value = 0.01063829787234043
$ws.Cells.Item(3, 3).Value = 0.02659574468085106
$ws.Cells.Item(3, 10).Value = 0.05851063829787234
$ws.Cells.Item(3, 16).Value = 0.675531914893617
$ws.Cells.Item(3, 19).Value = 0.2287234042553191
$ws.Cells.Item(4, 15).Value = 0.02857142857142857
$ws.Cells.Item(4, 16).Value = 0.5714285714285714
$ws.Cells.Item(4, 19).Value = 0.4
$ws.Cells.Item(6, 2).Value = 0.06696428571428571
$ws.Cells.Item(6, 4).Value = 0.004464285714285714
$ws.Cells.Item(6, 5).Value = 0.004464285714285714
$ws.Cells.Item(6, 6).Value = 0.0625
$ws.Cells.Item(6, 10).Value = 0.2321428571428572
$ws.Cells.Item(6, 15).Value = 0.04910714285714286
$ws.Cells.Item(6, 17).Value = 0.1785714285714286
$ws.Cells.Item(6, 18).Value = 0.07142857142857142
$ws.Cells.Item(6, 19).Value = 0.3303571428571428
$ws.Cells.Item(7, 2).Value = 0.1446280991735537
$ws.Cells.Item(7, 4).Value = 0.01652892561983471
$ws.Cells.Item(7, 6).Value = 0.04545454545454546
$ws.Cells.Item(7, 10).Value = 0.1363636363636364
$ws.Cells.Item(7, 15).Value = 0.008264462809917356
$ws.Cells.Item(7, 17).Value = 0.1652892561983471
$ws.Cells.Item(7, 18).Value = 0.07851239669421488
$ws.Cells.Item(7, 19).Value = 0.4049586776859504
$ws.Cells.Item(8, 2).Value = 0.0975103734439834
$ws.Cells.Item(8, 4).Value = 0.01452282157676349
$ws.Cells.Item(8, 5).Value = 0.002074688796680498
$ws.Cells.Item(8, 6).Value = 0.05601659751037345
$ws.Cells.Item(8, 10).Value = 0.1431535269709543
$ws.Cells.Item(8, 15).Value = 0.01659751037344398
$ws.Cells.Item(8, 17).Value = 0.1390041493775934
$ws.Cells.Item(8, 18).Value = 0.1099585062240664
$ws.Cells.Item(8, 19).Value = 0.4211618257261411
$ws.Cells.Item(9, 2).Value = 0.09036144578313253
$ws.Cells.Item(9, 4).Value = 0.006024096385542169
$ws.Cells.Item(9, 6).Value = 0.0783132530120482
$ws.Cells.Item(9, 10).Value = 0.1144578313253012
$ws.Cells.Item(9, 15).Value = 0.02409638554216868
$ws.Cells.Item(9, 17).Value = 0.1927710843373494
$ws.Cells.Item(9, 18).Value = 0.1204819277108434
$ws.Cells.Item(9, 19).Value = 0.3734939759036144
$ws.Cells.Item(10, 2).Value = 0.1177761521580102
$ws.Cells.Item(10, 4).Value = 0.01755669348939283
$ws.Cells.Item(10, 5).Value = 0.000731528895391368
$ws.Cells.Item(10, 6).Value = 0.06949524506217995
$ws.Cells.Item(10, 10).Value = 0.1382589612289686
$ws.Cells.Item(10, 15).Value = 0.01536210680321873
$ws.Cells.Item(10, 17).Value = 0.1799561082662765
$ws.Cells.Item(10, 18).Value = 0.08558888076079005
$ws.Cells.Item(10, 19).Value = 0.3752743233357718
$ws.Cells.Item(11, 7).Value = 0.1450777202072539
$ws.Cells.Item(11, 10).Value = 0.1010362694300518
$ws.Cells.Item(11, 11).Value = 0.2124352331606218
$ws.Cells.Item(11, 12).Value = 0.5336787564766839
$ws.Cells.Item(11, 19).Value = 0.007772020725388601
$ws.Cells.Item(12, 7).Value = 0.7285067873303167
$ws.Cells.Item(12, 10).Value = 0.1855203619909502
$ws.Cells.Item(12, 12).Value = 0.04072398190045249
$ws.Cells.Item(12, 19).Value = 0.04524886877828054
$ws.Cells.Item(13, 7).Value = 0.8108108108108109
$ws.Cells.Item(13, 10).Value = 0.1621621621621622
$ws.Cells.Item(13, 19).Value = 0.02702702702702703
$ws.Cells.Item(15, 6).Value = 0.03543307086614173
$ws.Cells.Item(15, 8).Value = 0.1889763779527559
$ws.Cells.Item(15, 9).Value = 0.05905511811023622
$ws.Cells.Item(15, 10).Value = 0.3622047244094488
$ws.Cells.Item(15, 11).Value = 0.06299212598425197
$ws.Cells.Item(15, 13).Value = 0.003937007874015748
$ws.Cells.Item(15, 14).Value = 0.003937007874015748
$ws.Cells.Item(15, 15).Value = 0.05511811023622047
$ws.Cells.Item(15, 19).Value = 0.2283464566929134
$ws.Cells.Item(16, 6).Value = 0.005235602094240838
$ws.Cells.Item(16, 8).Value = 0.1518324607329843
$ws.Cells.Item(16, 9).Value = 0.06282722513089005
$ws.Cells.Item(16, 10).Value = 0.450261780104712
$ws.Cells.Item(16, 11).Value = 0.08900523560209424
$ws.Cells.Item(16, 13).Value = 0.01570680628272251
$ws.Cells.Item(16, 14).Value = 0.005235602094240838
$ws.Cells.Item(16, 15).Value = 0.08900523560209424
$ws.Cells.Item(16, 19).Value = 0.1308900523560209
$ws.Cells.Item(17, 6).Value = 0.007159904534606206
$ws.Cells.Item(17, 8).Value = 0.1909307875894988
$ws.Cells.Item(17, 9).Value = 0.05250596658711217
$ws.Cells.Item(17, 10).Value = 0.441527446300716
$ws.Cells.Item(17, 11).Value = 0.1241050119331742
$ws.Cells.Item(17, 13).Value = 0.01909307875894988
$ws.Cells.Item(17, 15).Value = 0.081145584725537
$ws.Cells.Item(17, 19).Value = 0.08353221957040573
$ws.Cells.Item(18, 6).Value = 0.01333333333333333
$ws.Cells.Item(18, 8).Value = 0.1822222222222222
$ws.Cells.Item(18, 9).Value = 0.09333333333333334
$ws.Cells.Item(18, 10).Value = 0.4266666666666667
$ws.Cells.Item(18, 11).Value = 0.1022222222222222
$ws.Cells.Item(18, 13).Value = 0.02222222222222222
$ws.Cells.Item(18, 15).Value = 0.05333333333333334
$ws.Cells.Item(18, 19).Value = 0.1066666666666667
$ws.Cells.Item(19, 6).Value = 0.01665404996214989
$ws.Cells.Item(19, 8).Value = 0.2172596517789553
$ws.Cells.Item(19, 9).Value = 0.07267221801665405
$ws.Cells.Item(19, 10).Value = 0.3459500378501135
$ws.Cells.Item(19, 11).Value = 0.1438304314912945
$ws.Cells.Item(19, 13).Value = 0.0174110522331567
$ws.Cells.Item(19, 14).Value = 0.001514004542013626
$ws.Cells.Item(19, 15).Value = 0.07494322482967448
$ws.Cells.Item(19, 19).Value = 0.1097653292959879
